$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.234.86"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "1.644.77"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.59%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").Value = "1.873.74"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "1.646.92"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.545"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.83%  "

$ws.Range("D17").Value = "27.210.90"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("E18").Value = "  +1.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.93%  "

$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  +1.12%  "

$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("E34").Value = "  +1.22%  "

$ws.Range("D35").Value = "1.276.10"
$ws.Range("E35").Value = "  +2.36%  "

$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.860"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.10%  "

$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("E42").Value = "  +6.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").Value = "1.784.08"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("E47").Value = "  +1.91%  "

$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "
